$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: record count for Hierarchical OPT ---
$ws.Range("C3").Value = 10

# --- Row 5: NO (comment on naive hierarchical run) ---
$ws.Range("C5").Value = "NO"

# --- Row 6 (new): 4 means opt(ALB) ---
$ws.Range("A6").Value = "4 means opt(ALB)"
$ws.Range("C6").Value = 10
$ws.Range("E6").Value = "da calcolare il tempo( 1.5s)"

# --- Row 7 (new): btw parallel(ALB) ---
$ws.Range("A7").Value = "btw parallel(ALB)"
$ws.Range("C7").Value = 1
$ws.Range("E7").Value = "da calcolare il tempo(40 min)"

# --- Row 8 (new): btw nayve ---
$ws.Range("A8").Value = "btw nayve"
$ws.Range("B8").Value = "caricare pkl"
$ws.Range("C8").Value = "NO"

# --- Row 9 (new): 4 means Nayve ---
$ws.Range("A9").Value = "4 means Nayve"
$ws.Range("B9").Value = "caricare pkl"
$ws.Range("C9").Value = "NO"

# --- Row 10 (new): spectral 80% comp connessa ---
$ws.Range("A10").Value = "spectral 80% comp connessa"
$ws.Range("B10").Value = "caricare pkl"
$ws.Range("C10").Value = 1

# --- Row 11 (new): spectral 70% comp connessa ---
$ws.Range("A11").Value = "spectral 70% comp connessa"
$ws.Range("B11").Value = "caricare pkl"
$ws.Range("C11").Value = 1

# --- Row 12 (new): btw 50% comp connessa ---
$ws.Range("A12").Value = "btw 50% comp connessa"
$ws.Range("B12").Value = "caricare pkl"
$ws.Range("C12").Value = 1

# --- Row 13 (new): btw 70% comp connessa ---
$ws.Range("A13").Value = "btw 70% comp connessa"
$ws.Range("B13").Value = "caricare pkl"
$ws.Range("C13").Value = 1

# --- Row 17 (new): link pkl ---
$ws.Range("A17").Value = "link pkl"
$ws.Range("B17").Value = "https://drive.google.com/drive/folders/10g-AdWR3G9_fo3vCiGeoRuXzTAIylEBr?usp=sharing"

# --- Column widths widened slightly to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 31.0
$ws.Columns.Item(2).ColumnWidth = 32.666666666666664
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
$ws.Columns.Item(5).ColumnWidth = 25.0

# --- Final selection left on F19, matching the saved workbook state ---
$ws.Range("F19").Select()
